# Applies the "Actualizacion automatica 2025-10-14 08:30:10" update to
# /tmp/work/before.xlsx across its three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("E4").Value  = 64.79000000000001
$ws1.Range("M4").Value  = 2073.71

$ws1.Range("D12").Value = 475.2

$ws1.Range("D28").Value = 457.92
$ws1.Range("M28").Value = 109.45

$ws1.Range("D36").Value = 475.2
$ws1.Range("I36").Value = 778.91

# Row 55 "X de 53" counters
$ws1.Range("D55").Value = "6 de 53"
$ws1.Range("E55").Value = "3 de 53"
$ws1.Range("I55").Value = "7 de 53"
$ws1.Range("M55").Value = "7 de 53"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F4").Value  = 3693.7
$ws2.Range("F12").Value = 475.2
$ws2.Range("F28").Value = 567.37
$ws2.Range("F36").Value = 1254.11
$ws2.Range("F59").Value = 20419.1

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 4190.4
$ws3.Range("E3").Value = 13478.7470988183
$ws3.Range("F3").Value = 0.2371591552531843

$ws3.Range("D4").Value = 1319.38
$ws3.Range("E4").Value = -276.1571147347202
$ws3.Range("F4").Value = 1.264715353387303

$ws3.Range("D7").Value = 1535.81
$ws3.Range("E7").Value = -649.0989837124259
$ws3.Range("F7").Value = 1.732029908041554

$ws3.Range("D12").Value = 8253.01
$ws3.Range("E12").Value = 44410.11
$ws3.Range("F12").Value = 0.1567132748686367

$ws3.Range("D14").Value = 18398.41
$ws3.Range("E14").Value = 80618.09661190613
$ws3.Range("F14").Value = 0.1858115442520339

# Column E widened from stored width 23 to 24 (character units as read via
# the ColumnWidth COM property are offset by ~0.83 from the OOXML <col>
# width, so 23.17 here serializes back out to width="24").
$ws3.Columns.Item(5).ColumnWidth = 23.17
